$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("by_prov")
$ws2 = $wb.Worksheets.Item("all")

# --- sheet "by_prov": update OR / lower / higher figures with the new PSM model ---
$ws1.Range("B2").Value = 3.04
$ws1.Range("C2").Value = 1.84
$ws1.Range("D2").Value = 5.25

$ws1.Range("B3").Value = 1.34
$ws1.Range("C3").Value = 1.1399999999999999
$ws1.Range("D3").Value = 1.59

$ws1.Range("B4").Value = 3.46
$ws1.Range("C4").Value = 2.29
$ws1.Range("D4").Value = 5.35

$ws1.Range("B5").Value = 1.5
$ws1.Range("C5").Value = 1.27
$ws1.Range("D5").Value = 1.78

$ws1.Range("B6").Value = 4.03
$ws1.Range("C6").Value = 2.87
$ws1.Range("D6").Value = 5.88

$ws1.Range("B7").Value = 1.42
$ws1.Range("C7").Value = 1.23
$ws1.Range("D7").Value = 1.65

$ws1.Range("B8").Value = 1.95
$ws1.Range("C8").Value = 1.37
$ws1.Range("D8").Value = 2.8

$ws1.Range("B9").Value = 1.49
$ws1.Range("C9").Value = 1.27
$ws1.Range("D9").Value = 1.76

$ws1.Range("B10").Value = 7.56
$ws1.Range("C10").Value = 4.7300000000000004
$ws1.Range("D10").Value = 12.8

$ws1.Range("B11").Value = 1.45
$ws1.Range("C11").Value = 1.24
$ws1.Range("D11").Value = 1.7

# --- sheet "all": add new "Unweighted Analysis" block above the existing (now "weighted") table ---
$ws2.Rows("1:6").Insert()

$ws2.Range("A1").Value = "Unweighted Analysis"

$ws2.Range("A2").Value = "group"
$ws2.Range("B2").Value = "OR"
$ws2.Range("C2").Value = "lower"
$ws2.Range("D2").Value = "higher"

$ws2.Range("A3").Value = "Anti-N Positivity"
$ws2.Range("A4").Value = "Anti-S Positivity"

# --- update the (now shifted-down) weighted-analysis figures with the new PSM model ---
$ws2.Range("B8").Value = 3.23
$ws2.Range("C8").Value = 2.74
$ws2.Range("D8").Value = 3.83

$ws2.Range("B9").Value = 1.37
$ws2.Range("C9").Value = 1.28
$ws2.Range("D9").Value = 1.47

# --- view state: active sheet/tab + selection to match the final saved state ---
$ws2.Select() | Out-Null
$ws2.Range("D10").Select() | Out-Null

$ws1.Select() | Out-Null
$ws1.Range("D12").Select() | Out-Null
